$wb = $excel.ActiveWorkbook

# This script re-applies a market-data refresh (scheduled runner update) to the
# FFXIV Hades-server crafting-profit workbook: columns H-N on several sheets hold
# current market-board prices/profit calcs that get replaced wholesale on each run.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1664.8667
$ws.Range("I28").Value = 1652.091
$ws.Range("J28").Value = 1700
$ws.Range("K28").Value = 1652.091
$ws.Range("L28").Value = 1700
$ws.Range("M28").Value = -1167.091
$ws.Range("N28").Value = -2670
$ws.Range("H62").Value = 3937.5
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 4583.3335
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 4583.3335
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -5831.3335
$ws.Range("H65").Value = 3937.5
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 4583.3335
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 22916.6675
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -29156.6675
$ws.Range("H74").Value = 3925
$ws.Range("I74").Value = 3933.3333
$ws.Range("J74").Value = 3920
$ws.Range("K74").Value = 3933.3333
$ws.Range("L74").Value = 3920
$ws.Range("M74").Value = -2997.3333
$ws.Range("N74").Value = -5792
$ws.Range("H77").Value = 3925
$ws.Range("I77").Value = 3933.3333
$ws.Range("J77").Value = 3920
$ws.Range("K77").Value = 19666.6665
$ws.Range("L77").Value = 19600
$ws.Range("M77").Value = -14986.6665
$ws.Range("N77").Value = -28960
$ws.Range("H113").Value = 3452.9412
$ws.Range("I113").Value = 3169.2307
$ws.Range("J113").Value = 4375
$ws.Range("K113").Value = 3169.2307
$ws.Range("L113").Value = 4375
$ws.Range("M113").Value = 84.76929999999993
$ws.Range("N113").Value = -10883
$ws.Range("H132").Value = 805223.1
$ws.Range("I132").Value = 1705.5098
$ws.Range("J132").Value = 4903163
$ws.Range("K132").Value = 5116.5294
$ws.Range("L132").Value = 14709489
$ws.Range("M132").Value = -2586.5294
$ws.Range("N132").Value = -14714549

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14740.6045
$ws.Range("I32").Value = 17754.436
$ws.Range("J32").Value = 7038.593
$ws.Range("K32").Value = 17754.436
$ws.Range("L32").Value = 7038.593
$ws.Range("M32").Value = -17467.436
$ws.Range("N32").Value = -7612.593
$ws.Range("H45").Value = 1305.7778
$ws.Range("I45").Value = 921.7143
$ws.Range("J45").Value = 2650
$ws.Range("K45").Value = 921.7143
$ws.Range("L45").Value = 2650
$ws.Range("M45").Value = -544.7143
$ws.Range("N45").Value = -3404
$ws.Range("H61").Value = 250501000
$ws.Range("I61").Value = 500500000
$ws.Range("K61").Value = 500500000
$ws.Range("M61").Value = -500499788
$ws.Range("H74").Value = 6812465
$ws.Range("I74").Value = 11954064
$ws.Range("J74").Value = 64115.5
$ws.Range("K74").Value = 11954064
$ws.Range("L74").Value = 64115.5
$ws.Range("M74").Value = -11953190
$ws.Range("N74").Value = -65863.5
$ws.Range("H77").Value = 6812465
$ws.Range("I77").Value = 11954064
$ws.Range("J77").Value = 64115.5
$ws.Range("K77").Value = 59770320
$ws.Range("L77").Value = 320577.5
$ws.Range("M77").Value = -59765952
$ws.Range("N77").Value = -329313.5
$ws.Range("H97").Value = 2718265.8
$ws.Range("I97").Value = 4167507.2
$ws.Range("K97").Value = 4167507.2
$ws.Range("M97").Value = -4167011.2
$ws.Range("H102").Value = 9532223
$ws.Range("I102").Value = 10998350
$ws.Range("J102").Value = 2400
$ws.Range("K102").Value = 10998350
$ws.Range("L102").Value = 2400
$ws.Range("M102").Value = -10996728
$ws.Range("N102").Value = -5644
$ws.Range("H113").Value = 30000
$ws.Range("J113").Value = 30000
$ws.Range("L113").Value = 30000
$ws.Range("N113").Value = -38678
$ws.Range("H122").Value = 27780400
$ws.Range("I122").Value = 3496.6667
$ws.Range("K122").Value = 10490.0001
$ws.Range("M122").Value = -8040.000100000001
$ws.Range("H132").Value = 191155.45
$ws.Range("I132").Value = 144085.14
$ws.Range("J132").Value = 273528.5
$ws.Range("K132").Value = 432255.42
$ws.Range("L132").Value = 820585.5
$ws.Range("M132").Value = -429725.42
$ws.Range("N132").Value = -825645.5
$ws.Range("H136").Value = 250501000
$ws.Range("I136").Value = 500500000
$ws.Range("K136").Value = 1501500000
$ws.Range("M136").Value = -1501497450

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 28361.092
$ws.Range("J82").Value = 40281.715
$ws.Range("L82").Value = 40281.715
$ws.Range("N82").Value = -41047.715
$ws.Range("H85").Value = 28361.092
$ws.Range("J85").Value = 40281.715
$ws.Range("L85").Value = 40281.715
$ws.Range("N85").Value = -42933.715
$ws.Range("H94").Value = 518.38464
$ws.Range("I94").Value = 319.875
$ws.Range("K94").Value = 319.875
$ws.Range("M94").Value = 131.125
$ws.Range("H134").Value = 2948
$ws.Range("I134").Value = 1638.9286
$ws.Range("K134").Value = 4916.7858
$ws.Range("M134").Value = -2381.7858

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 111111630
$ws.Range("I22").Value = 200000240
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 200000240
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -199999890
$ws.Range("N22").Value = -1575
$ws.Range("H31").Value = 3939.6206
$ws.Range("I31").Value = 1469.3889
$ws.Range("J31").Value = 7981.8184
$ws.Range("K31").Value = 1469.3889
$ws.Range("L31").Value = 7981.8184
$ws.Range("M31").Value = -1174.3889
$ws.Range("N31").Value = -8571.8184
$ws.Range("H34").Value = 3939.6206
$ws.Range("I34").Value = 1469.3889
$ws.Range("J34").Value = 7981.8184
$ws.Range("K34").Value = 1469.3889
$ws.Range("L34").Value = 7981.8184
$ws.Range("M34").Value = -1267.3889
$ws.Range("N34").Value = -8385.8184
$ws.Range("H41").Value = 10500
$ws.Range("H50").Value = 24873.6
$ws.Range("J50").Value = 24873.6
$ws.Range("L50").Value = 24873.6
$ws.Range("N50").Value = -26123.6
$ws.Range("H51").Value = 28099
$ws.Range("J51").Value = 28099
$ws.Range("L51").Value = 28099
$ws.Range("N51").Value = -29571
$ws.Range("H60").Value = 14071.429
$ws.Range("J60").Value = 14071.429
$ws.Range("L60").Value = 14071.429
$ws.Range("N60").Value = -15093.429
$ws.Range("H61").Value = 28099
$ws.Range("J61").Value = 28099
$ws.Range("L61").Value = 28099
$ws.Range("N61").Value = -28795
$ws.Range("H122").Value = 2414
$ws.Range("I122").Value = 1902
$ws.Range("J122").Value = 3438
$ws.Range("K122").Value = 5706
$ws.Range("L122").Value = 10314
$ws.Range("M122").Value = -3256
$ws.Range("N122").Value = -15214

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1864.6875
$ws.Range("I109").Value = 667.125
$ws.Range("J109").Value = 3062.25
$ws.Range("K109").Value = 2001.375
$ws.Range("L109").Value = 9186.75
$ws.Range("M109").Value = -961.375
$ws.Range("N109").Value = -11266.75
$ws.Range("H114").Value = 33333936
$ws.Range("I114").Value = 622.125
$ws.Range("J114").Value = 166667180
$ws.Range("K114").Value = 1866.375
$ws.Range("L114").Value = 500001540
$ws.Range("M114").Value = 1387.625
$ws.Range("N114").Value = -500008048
$ws.Range("H117").Value = 3704712
$ws.Range("I117").Value = 495.22223
$ws.Range("J117").Value = 7408929
$ws.Range("K117").Value = 1485.66669
$ws.Range("L117").Value = 22226787
$ws.Range("M117").Value = 1956.33331
$ws.Range("N117").Value = -22233671

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 39850
$ws.Range("J103").Value = 39850
$ws.Range("L103").Value = 39850
$ws.Range("N103").Value = -42194
$ws.Range("H122").Value = 4722
$ws.Range("I122").Value = 4364
$ws.Range("J122").Value = 5080
$ws.Range("K122").Value = 13092
$ws.Range("L122").Value = 15240
$ws.Range("M122").Value = -10642
$ws.Range("N122").Value = -20140
$ws.Range("H126").Value = 2390.8333
$ws.Range("I126").Value = 1831.6666
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 5494.9998
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -3024.9998
$ws.Range("N126").Value = -13790
$ws.Range("H129").Value = 48228.184
$ws.Range("I129").Value = 45890
$ws.Range("J129").Value = 49105
$ws.Range("K129").Value = 45890
$ws.Range("L129").Value = 49105
$ws.Range("M129").Value = -40890
$ws.Range("N129").Value = -59105
$ws.Range("H132").Value = 113053.055
$ws.Range("I132").Value = 73357.5
$ws.Range("K132").Value = 220072.5
$ws.Range("M132").Value = -217542.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1390.5714
$ws.Range("I100").Value = 1200.125
$ws.Range("K100").Value = 1200.125
$ws.Range("M100").Value = -659.125
$ws.Range("H119").Value = 20420
$ws.Range("J119").Value = 20420
$ws.Range("L119").Value = 20420
$ws.Range("N119").Value = -30096
$ws.Range("H122").Value = 3157.1072
$ws.Range("I122").Value = 2550.5715
$ws.Range("K122").Value = 7651.7145
$ws.Range("M122").Value = -5201.7145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2559.8333
$ws.Range("I122").Value = 2289
$ws.Range("J122").Value = 2999.9375
$ws.Range("K122").Value = 6867
$ws.Range("L122").Value = 8999.8125
$ws.Range("M122").Value = -4417
$ws.Range("N122").Value = -13899.8125
$ws.Range("H136").Value = 169792
$ws.Range("I136").Value = 201500.8
$ws.Range("K136").Value = 604502.3999999999
$ws.Range("M136").Value = -601952.3999999999
